# Apply the edits described by the diff:
# 1. Change the shared string used by A10 to "Should not be here"
# 2. Change font color (2nd font, theme=1) to rgb FF000000 -- applies to column B cells
# 3. Change row height for rows 1-10 from 17.25 to 19.5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the text value of A10 (which uses the shared string replaced in the diff)
$ws.Range("A10").Value = "Should not be here"

# 2. Update font color for column B (which uses the second font / style s="2")
$ws.Range("B1:B10").Font.Color = 0

# 3. Update row heights for rows 1 through 10
for ($r = 1; $r -le 10; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}
